$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45044

# Row 3
$ws.Range("D3").Value = 45044
$ws.Range("M3").Value = 40

# Row 4
$ws.Range("D4").Value = 45049
$ws.Range("M4").Value = 50

# Row 5
$ws.Range("D5").Value = 45049
$ws.Range("M5").Value = 60

# Row 6
$ws.Range("D6").Value = 44699
$ws.Range("M6").Value = 60
$ws.Range("Q6").Value = '$/caja 15 kilos granel'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("S6").Value = 867
$ws.Range("T6").Value = 15

# Row 7
$ws.Range("D7").Value = 44699
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 11000
$ws.Range("P7").Value = 11500
$ws.Range("Q7").Value = '$/caja 15 kilos granel'
$ws.Range("R7").Value = 'Provincia de Curicó'
$ws.Range("S7").Value = 767
$ws.Range("T7").Value = 15

# Row 8
$ws.Range("D8").Value = 45033
$ws.Range("L8").Value = 'Especial'
$ws.Range("N8").Value = 13000
$ws.Range("O8").Value = 13000
$ws.Range("P8").Value = 13000
$ws.Range("Q8").Value = '$/caja 18 kilos empedrada'
$ws.Range("S8").Value = 722

# Row 9
$ws.Range("D9").Value = 45033
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("S9").Value = 667

# Row 10
$ws.Range("D10").Value = 45021
$ws.Range("M10").Value = 50
$ws.Range("Q10").Value = '$/caja 18 kilos granel'

# Row 11
$ws.Range("D11").Value = 45050
$ws.Range("L11").Value = 'Especial'
$ws.Range("N11").Value = 13000
$ws.Range("O11").Value = 13000
$ws.Range("P11").Value = 13000
$ws.Range("Q11").Value = '$/caja 18 kilos empedrada'
$ws.Range("S11").Value = 722

# Row 12
$ws.Range("D12").Value = 45050
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range("Q12").Value = '$/caja 18 kilos empedrada'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 667
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 45043
$ws.Range("L13").Value = 'Especial'
$ws.Range("M13").Value = 40
$ws.Range("N13").Value = 13000
$ws.Range("O13").Value = 13000
$ws.Range("P13").Value = 13000
$ws.Range("Q13").Value = '$/caja 18 kilos empedrada'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 722
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("D14").Value = 45043
$ws.Range("L14").Value = 'Primera'
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("S14").Value = 667

# Row 15
$ws.Range("D15").Value = 45020
$ws.Range("Q15").Value = '$/caja 18 kilos granel'

